$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.406.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.93%  '
$ws.Range("D3").Value = "'3.770.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.90%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = "'605.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.64%  '
$ws.Range("D6").Value = "'170.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.81%  '
$ws.Range("D7").Value = "'3.767.90"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.99%  '
$ws.Range("D9").Value = "'0.540"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.18%  '
$ws.Range("E10").Value = '  +5.51%  '
$ws.Range("E11").Value = '  +3.19%  '
$ws.Range("E12").Value = '  +0.54%  '
$ws.Range("D13").Value = "'38.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.51%  '
$ws.Range("D14").Value = "'0.0000252"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.62%  '
$ws.Range("D15").Value = "'4.397.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.84%  '
$ws.Range("D16").Value = "'3.759.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.60%  '
$ws.Range("D17").Value = "'69.394.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.86%  '
$ws.Range("D18").Value = "'7.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.34%  '
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("D21").Value = "'10.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +19.33%  '
$ws.Range("D22").Value = "'497.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.29%  '
$ws.Range("D23").Value = "'0.732"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.94%  '
$ws.Range("D24").Value = "'0.0000156"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +12.42%  '
$ws.Range("D25").Value = "'85.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'2.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.11%  '
$ws.Range("D27").Value = "'12.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.00%  '
$ws.Range("D28").Value = "'10.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.95%  '
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("E30").Value = '  +7.77%  '
$ws.Range("D31").Value = "'3.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.25%  '
$ws.Range("D32").Value = "'7.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.65%  '
$ws.Range("D33").Value = "'32.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.83%  '
$ws.Range("D34").Value = "'3.914.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.96%  '
$ws.Range("E35").Value = '  +1.60%  '
$ws.Range("D36").Value = "'3.704.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.73%  '
$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("E38").Value = '  +2.45%  '
$ws.Range("D39").Value = "'5.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.96%  '
$ws.Range("E40").Value = '  +1.96%  '
$ws.Range("E41").Value = '  +1.10%  '
$ws.Range("D42").Value = "'3.03"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +8.92%  '
$ws.Range("D43").Value = "'439.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.37%  '
$ws.Range("D44").Value = "'48.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("D45").Value = "'1.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.42%  '
$ws.Range("D46").Value = "'8.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.59%  '
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("D48").Value = "'40.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("D49").Value = "'2.820.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.39%  '
$ws.Range("D50").Value = "'141.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.06%  '
$ws.Range("D51").Value = "'0.0357"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.76%  '
